$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage for numeric-looking
# strings, so Excel does not auto-convert them to numbers (the source data
# uses plain text for these columns, e.g. "1.003", "310.73", percentages).
function Set-TextValue($sheet, $cellRef, [string]$value) {
    $cell = $sheet.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "27.780.55"
Set-TextValue $ws "E2" "  -2.66%  "
Set-TextValue $ws "D3" "1.780.96"
Set-TextValue $ws "E3" "  -2.18%  "
Set-TextValue $ws "E4" "  +0.11%  "
Set-TextValue $ws "D5" "310.73"
Set-TextValue $ws "E5" "  -1.77%  "
Set-TextValue $ws "E6" "  +0.12%  "
Set-TextValue $ws "D7" "0.5126"
Set-TextValue $ws "E7" "  -0.57%  "
Set-TextValue $ws "D8" "0.3779"
Set-TextValue $ws "E8" "  -2.41%  "
Set-TextValue $ws "D9" "0.07766"
Set-TextValue $ws "E9" "  -7.77%  "
Set-TextValue $ws "D10" "41.13"
Set-TextValue $ws "E10" "  -1.84%  "
Set-TextValue $ws "D11" "1.084"
Set-TextValue $ws "E11" "  -2.37%  "
Set-TextValue $ws "E12" "  +0.14%  "
Set-TextValue $ws "D13" "6.192"
Set-TextValue $ws "E13" "  -3.21%  "
Set-TextValue $ws "D14" "20.09"
Set-TextValue $ws "E14" "  -4.51%  "
Set-TextValue $ws "D15" "1.773.05"
Set-TextValue $ws "E15" "  -2.19%  "
Set-TextValue $ws "D16" "7.160"
Set-TextValue $ws "E16" "  -4.61%  "
Set-TextValue $ws "D17" "91.65"
Set-TextValue $ws "E17" "  -1.33%  "
Set-TextValue $ws "D18" "0.00001069"
Set-TextValue $ws "E18" "  -5.65%  "
Set-TextValue $ws "D19" "0.06532"
Set-TextValue $ws "E19" "  -2.36%  "
Set-TextValue $ws "E21" "  -4.50%  "
Set-TextValue $ws "D22" "5.909"
Set-TextValue $ws "E22" "  -2.93%  "
Set-TextValue $ws "D23" "27.832.99"
Set-TextValue $ws "E23" "  -2.57%  "
Set-TextValue $ws "D24" "10.95"
Set-TextValue $ws "E24" "  -4.19%  "
Set-TextValue $ws "D25" "2.240"
Set-TextValue $ws "E25" "  -1.77%  "
Set-TextValue $ws "D26" "158.60"
Set-TextValue $ws "E26" "  -0.20%  "
Set-TextValue $ws "D27" "20.17"
Set-TextValue $ws "E27" "  -4.37%  "
Set-TextValue $ws "D28" "1.982.98"
Set-TextValue $ws "E28" "  -2.01%  "
Set-TextValue $ws "D29" "2.345"
Set-TextValue $ws "E29" "  -3.03%  "
Set-TextValue $ws "D30" "125.10"
Set-TextValue $ws "E30" "  -0.78%  "
Set-TextValue $ws "D31" "0.1076"
Set-TextValue $ws "E31" "  -0.34%  "
Set-TextValue $ws "D32" "1.025"
Set-TextValue $ws "E32" "  -6.40%  "
Set-TextValue $ws "E33" "  -1.80%  "
Set-TextValue $ws "D34" "5.473"
Set-TextValue $ws "E34" "  -4.79%  "
Set-TextValue $ws "D35" "0.07049"
Set-TextValue $ws "E35" "  -7.24%  "
Set-TextValue $ws "E36" "  -2.64%  "
Set-TextValue $ws "D37" "8.688"
Set-TextValue $ws "E37" "  -0.53%  "
Set-TextValue $ws "D38" "0.2116"
Set-TextValue $ws "E38" "  -4.91%  "
Set-TextValue $ws "E39" "  +2.32%  "
Set-TextValue $ws "D40" "5.006"
Set-TextValue $ws "E40" "  -3.83%  "
Set-TextValue $ws "D41" "0.6073"
Set-TextValue $ws "E41" "  -4.13%  "
Set-TextValue $ws "E42" "  +0.02%  "
Set-TextValue $ws "D43" "1.150"
Set-TextValue $ws "E43" "  -3.23%  "
Set-TextValue $ws "D44" "1.319"
Set-TextValue $ws "E44" "  -5.81%  "
Set-TextValue $ws "B45" "EnergySwap"
Set-TextValue $ws "C45" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws "D45" "13.09"
Set-TextValue $ws "E45" "  -3.26%  "
Set-TextValue $ws "B46" "Decentraland"
Set-TextValue $ws "C46" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws "D46" "0.5948"
Set-TextValue $ws "E46" "  +0.47%  "
Set-TextValue $ws "D47" "3.714"
Set-TextValue $ws "E47" "  -1.37%  "
Set-TextValue $ws "D48" "127.26"
Set-TextValue $ws "E48" "  +1.04%  "
Set-TextValue $ws "D49" "1.211"
Set-TextValue $ws "E49" "  +0.93%  "
Set-TextValue $ws "D50" "1.893"
Set-TextValue $ws "E50" "  -4.86%  "
Set-TextValue $ws "D51" "0.06702"
Set-TextValue $ws "E51" "  -4.05%  "
